# Add a new "Italy" market test-data sheet, cloned from the existing
# "Slovakia" sheet (same layout/styles), with the Italy-specific market
# name and product code filled in. Mirrors the selection/active-tab state
# left behind by the original author: Slovakia is no longer the active
# tab (its sheet now shows a "select all" selection), and Italy becomes
# the new active tab with cell G24 selected.

$wb = $excel.ActiveWorkbook

# Duplicate the last sheet (Slovakia) — this copies all formatting,
# merged cells, column widths, etc. — and drop the copy right after it.
$slovakia = $wb.Worksheets.Item("Slovakia")
$slovakia.Copy($null, $slovakia)

$italy = $wb.Worksheets.Item($slovakia.Index + 1)
$italy.Name = "Italy"

# Fill in the Italy-specific market name and product/ticket code.
$italy.Range("B2").Value = "Italy Market"
$italy.Range("B4").Value = "NGC-3145/T2226/T2447"

# Leave Slovakia with a "select everything" selection and no longer the
# active tab.
$slovakia.Select()
$slovakia.Cells.Select()

# Make Italy the active tab, with G24 as the selected cell.
$italy.Select()
$italy.Range("G24").Select()
